# SFD inputs workbook update
# - Vehicle Sections: Nosecone length corrected, selection left on C12
# - Aerodynamic Properties: proper inputs for "Off the rail" / "Max Q" rows
#   (acceleration/velocity/mach derived from real numbers), selection left on E2
# - A (now-empty) threaded comment/person is added and cleared, leaving the
#   person-list infrastructure part behind, matching the authored commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Aerodynamic Properties - fill in the proper/available inputs
# ---------------------------------------------------------------------
$aero = $wb.Worksheets.Item("Aerodynamic Properties")

# Row 2 = "Off the rail"
$aero.Range("C2").Formula = "=6.8*9.81"
$aero.Range("D2").Formula = "=29.7/343"

# Row 3 = "Max Q"
$aero.Range("B3").Formula = "=0.45*343"
$aero.Range("C3").Formula = "=7.47*9.81"
$aero.Range("D3").Value = 0.45

$aero.Range("E2").Select()

# ---------------------------------------------------------------------
# Comments infrastructure: a threaded comment was added while reviewing
# the sheet and then removed again, leaving the (now-empty) persons part
# behind in the package.
# ---------------------------------------------------------------------
$comment = $aero.Range("A1").AddCommentThreaded("Reviewed inputs")
$comment.Delete()

# ---------------------------------------------------------------------
# Sheet: Vehicle Sections (left active/selected on save)
# ---------------------------------------------------------------------
$vehicleSections = $wb.Worksheets.Item("Vehicle Sections")

# Nosecone Length (ft): 0.5 -> 1
$vehicleSections.Range("C2").Value = 1

# Leave the cursor where the author left it when they saved, with this
# sheet as the active tab
$vehicleSections.Activate()
$vehicleSections.Range("C12").Select()

$wb.Save()
